$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 15:22"

# Update row 6 data values (Galicia)
$ws.Range("B6").Value = 6946
$ws.Range("C6").Value = 997
$ws.Range("D6").Value = 5631
$ws.Range("E6").Value = 318
